# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (rId1), 演出 (rId2), 全部类型 (rId4) as per upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 124
$ws1.Range("F3").Value = 78
$ws1.Range("F4").Value = 718
$ws1.Range("F6").Value = 44
$ws1.Range("F7").Value = 2756
$ws1.Range("F9").Value = 1752
$ws1.Range("F12").Value = 714
$ws1.Range("F13").Value = 878
$ws1.Range("F14").Value = 154
$ws1.Range("F16").Value = 1110
$ws1.Range("F18").Value = 47
$ws1.Range("F20").Value = 6348
$ws1.Range("F21").Value = 247
$ws1.Range("F22").Value = 1381
$ws1.Range("F23").Value = 140
$ws1.Range("F26").Value = 296
$ws1.Range("F27").Value = 250
$ws1.Range("F28").Value = 60
$ws1.Range("F29").Value = 1091
$ws1.Range("F30").Value = 891
$ws1.Range("F32").Value = 84
$ws1.Range("F34").Value = 458
$ws1.Range("F35").Value = 1324
$ws1.Range("F37").Value = 138
$ws1.Range("F38").Value = 212
$ws1.Range("F39").Value = 3
$ws1.Range("F40").Value = 138
$ws1.Range("F41").Value = 178
$ws1.Range("F42").Value = 147

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 124
$ws4.Range("F3").Value = 78
$ws4.Range("F4").Value = 718
$ws4.Range("F6").Value = 15
$ws4.Range("F9").Value = 44
$ws4.Range("F10").Value = 2756
$ws4.Range("F12").Value = 1752
$ws4.Range("F15").Value = 714
$ws4.Range("F17").Value = 878
$ws4.Range("F18").Value = 154
$ws4.Range("F20").Value = 1110
$ws4.Range("F21").Value = 47
$ws4.Range("F23").Value = 6348
$ws4.Range("F24").Value = 247
$ws4.Range("F25").Value = 1381
$ws4.Range("F27").Value = 140
$ws4.Range("F30").Value = 296
$ws4.Range("F31").Value = 250
$ws4.Range("F32").Value = 60
$ws4.Range("F33").Value = 1091
$ws4.Range("F34").Value = 891
$ws4.Range("F36").Value = 84
$ws4.Range("F38").Value = 458
$ws4.Range("F39").Value = 1324
$ws4.Range("F41").Value = 138
$ws4.Range("F42").Value = 212
$ws4.Range("F43").Value = 3
$ws4.Range("F44").Value = 138
$ws4.Range("F45").Value = 178
$ws4.Range("F49").Value = 147
